$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 163, shifting the existing rows 163-178
# down to 165-180 (matches the target dimension A1:R180).
$ws.Range("A163:A164").EntireRow.Insert()

# New row 163: Ají / Cacho cabra rojo, Terminal Hortofrutícola Agro Chillán
$ws.Range("A163").Value = 7
$ws.Range("B163").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C163").Value = "Ñuble"
$ws.Range("D163").Value = 45021
$ws.Range("E163").Value = 16
$ws.Range("F163").Value = 100112021
$ws.Range("G163").Value = "Ají"
$ws.Range("H163").Value = "Cacho cabra rojo"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 50
$ws.Range("K163").Value = 15000
$ws.Range("L163").Value = 15000
$ws.Range("M163").Value = 15000
$ws.Range("N163").Value = "$/saco 25 kilos"
$ws.Range("O163").Value = "Región del Maule"
$ws.Range("P163").Value = 600
$ws.Range("Q163").Value = 25
$ws.Range("R163").Value = "Hortaliza"

# New row 164: Ají / Cristal, Terminal Hortofrutícola Agro Chillán
$ws.Range("A164").Value = 7
$ws.Range("B164").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C164").Value = "Ñuble"
$ws.Range("D164").Value = 45021
$ws.Range("E164").Value = 16
$ws.Range("F164").Value = 100112021
$ws.Range("G164").Value = "Ají"
$ws.Range("H164").Value = "Cristal"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 50
$ws.Range("K164").Value = 15000
$ws.Range("L164").Value = 15000
$ws.Range("M164").Value = 15000
$ws.Range("N164").Value = "$/saco 25 kilos"
$ws.Range("O164").Value = "Región del Maule"
$ws.Range("P164").Value = 600
$ws.Range("Q164").Value = 25
$ws.Range("R164").Value = "Hortaliza"
